$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# LogBook: the "Ideal/Actual remaining efforts" trailing columns (N:Q) for
# rows 15-20 burn down to 0 instead of staying at 2 (the shared-formula
# projection is overwritten with the final, actual value).
$ws.Range("N15:Q20").Value = 0

# Reflect the last-clicked cell from the edit session (was M20).
$ws.Range("N19").Select()

